$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily profit row (row 19) below the existing data.
# Force column A to be treated as plain text so the date string
# "09/05/2025" isn't auto-converted into a date serial number,
# matching the inline-string date cells already used in A2:A18.
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "09/05/2025"
# Reset the cell style back to Normal/General so no extra number
# format is left applied to the cell (matches the unstyled cells
# used by the rest of the date column).
$ws.Range("A19").Style = "Normal"

$ws.Range("B19").Value = 13761.69
